$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update I18 from a numeric 10 to a manual-status text value "[10]"
$ws.Range("I18").Value = "[10]"

# Set column widths: F (manualAudit/fastqFileName-ish wide text col) and G wider,
# to better display the long text values, matching the "manual status column" update.
$ws.Columns.Item(6).ColumnWidth = 73.3
$ws.Columns.Item(7).ColumnWidth = 68.1

# Move the active selection to I18 (the cell that was just edited)
$ws.Range("I18").Select() | Out-Null
